# Task 2.3: Triples for external and internal URIs
#
# Cleans up the "State Name" column in the zip code ranges table: several
# rows used ad-hoc / one-off labels (e.g. city call-outs or stray notes)
# instead of the canonical state name used elsewhere in the sheet. Collapse
# those variants down to the single canonical name so that the same state
# is referenced by one consistent string throughout (supports building
# clean triples keyed off the state name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old (one-off) State Name -> canonical State Name.
# NOTE: Missouri, New York State and District of Columbia are written first
# so that the brand-new shared strings are appended to the shared string
# table in that order.
$ws.Cells.Item(36, 2).Value2 = "Missouri"          # kc96 DataMO            -> Missouri

$ws.Cells.Item(48, 2).Value2 = "New York State"    # New York (Fishers Is)  -> New York State
$ws.Cells.Item(49, 2).Value2 = "New York State"    # New York               -> New York State

$ws.Cells.Item(11, 2).Value2 = "District of Columbia"  # Dist of Columbia   -> District of Columbia
$ws.Cells.Item(12, 2).Value2 = "District of Columbia"  # Dist of Columbia   -> District of Columbia
$ws.Cells.Item(13, 2).Value2 = "District of Columbia"  # Dist of Columbia   -> District of Columbia

$ws.Cells.Item(5, 2).Value2 = "Arkansas"           # Arkansas (Texarkana)   -> Arkansas
$ws.Cells.Item(17, 2).Value2 = "Georgia"           # Georga (Atlanta)       -> Georgia
$ws.Cells.Item(20, 2).Value2 = "Iowa"              # Iowa (OMAHA)           -> Iowa
$ws.Cells.Item(29, 2).Value2 = "Massachusetts"     # Massachusetts (Andover)-> Massachusetts
$ws.Cells.Item(38, 2).Value2 = "Mississippi"       # Mississippi(Warren)    -> Mississippi
$ws.Cells.Item(60, 2).Value2 = "Texas"             # Texas (Austin)         -> Texas
$ws.Cells.Item(63, 2).Value2 = "Texas"             # Texas (El Paso)        -> Texas

# Column B ("State Name") now holds shorter, uniform text - size the column
# to fit its (now narrower) contents.
$ws.Columns("B").ColumnWidth = 20

# Leave the selection where the edits were focused.
$ws.Range("B14").Select() | Out-Null
